$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting existing columns C.. to the right
$ws.Columns("C:C").Insert()

# Set header for the new column
$ws.Range("C1").Value = "pid"

# Fill in the pid values for rows 2-13
$pidValues = @(15,16,17,18,19,20,21,22,23,24,25,26)
for ($i = 0; $i -lt $pidValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $pidValues[$i]
}

# Update the view selection (also clears the stale topLeftCell scroll position)
$ws.Range("D17").Select()
